$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the bias/dark file labels & paths (rows 25-26), keep the flat label as-is
# but refresh its path value (row 27).
$ws.Range("A25").Value = "文件：偏置场"
$ws.Range("B25").Value = "/astro/294mm-pro/bias-master.fits"

$ws.Range("A26").Value = "文件：暗场"
$ws.Range("B26").Value = "/astro/294mm-pro/dark-master.fits"

$ws.Range("A27").Value = "文件目录：平场"
$ws.Range("B27").Value = "/astro/294mm-pro/flat/"

# Reflect the final cell selection left by the edit session.
$ws.Activate()
$ws.Range("B27").Select()
